# "Removed Type from imports"
#
# The "Type" column (G) - the Pool/CoInvest classifier - is no longer part
# of the portfolio-investments import template, so the whole column is
# deleted. Everything to the right (Folio No, Instrument, Currency,
# Investment Domicile) shifts one column to the left, and the data
# validations / cell comments that lived on those header cells need to
# move with them.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Grab the comment text that needs to "ride along" with the cells that are
# about to shift left, before the column structure changes underneath us.
$folioCommentText = $ws.Range("H1").Comment.Text()
$domicileCommentText = $ws.Range("K1").Comment.Text()

# Delete the entire "Type" column (G). Excel shifts everything after it
# one column to the left (H->G, I->H, J->I, K->J, ...).
$ws.Range("G:G").Delete()

# The comment that used to describe the old G column ("Pool / CoInvest")
# doesn't apply anymore, so replace it with the Folio No comment that used
# to sit on the old H1 (now G1).
$null = $ws.Range("G1").Comment.Text($folioCommentText)
$ws.Range("H1").Comment.Delete()

# The Investment Domicile comment used to live on K1; after the shift that
# column is now J1 (and it no longer has a comment of its own).
$ws.Range("K1").Comment.Delete()
$null = $ws.Range("J1").AddComment($domicileCommentText)

$ws.Range("G:G").Select()
